$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 22 and row 23 and need to be swapped.
$cols = @("A", "B", "E", "F", "G", "M", "Q", "R", "AC", "AM", "AO")

foreach ($col in $cols) {
    $addr22 = "$col" + "22"
    $addr23 = "$col" + "23"
    $val22 = $ws.Range($addr22).Value2
    $val23 = $ws.Range($addr23).Value2
    $ws.Range($addr22).Value2 = $val23
    $ws.Range($addr23).Value2 = $val22
}
